# Optimized regression run code, return time.
# Clear the RecvOrdStatus (J) and RecvErrorCode (K) columns for the sample
# rows in the report sheet. Rows that become fully empty as a result
# (rows 3 and 6, which only ever held J/K values) disappear from the
# worksheet entirely, while rows 4 and 5 keep their remaining Remark (M)
# values under their original row numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""

$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 11).Value = ""

$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = ""

$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = ""

$ws.Cells.Item(6, 10).Value = ""
$ws.Cells.Item(6, 11).Value = ""
